$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.979.58"
$ws.Range("E2").Value = "  +1.56%  "

$ws.Range("D3").Value = "3.421.96"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.72"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.96%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.60"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").Value = "4.007.34"
$ws.Range("E12").Value = "  +0.99%  "

$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.27"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.50%  "

$ws.Range("D15").Value = "3.413.67"
$ws.Range("E15").Value = "  +0.60%  "

$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").Value = "61.982.37"
$ws.Range("E17").Value = "  +1.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.99"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.71"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.43"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.553"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000116"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.193"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("E29").Value = "  +0.74%  "

$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("E31").Value = "  +2.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.56"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.28"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.81%  "

$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "168.06"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("D37").Value = "3.455.44"
$ws.Range("E37").Value = "  +0.95%  "

$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.48"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0754"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.787"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.45"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.68"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.18"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.47%  "

$ws.Range("D45").Value = "2.531.23"
$ws.Range("E45").Value = "  +3.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.89"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.45%  "

$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.62"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.88%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0264"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("E50").Value = "  -2.88%  "

$ws.Range("E51").Value = "  -0.44%  "
